# Cadastro.xlsx refactor — mirrors the "ExcelData / Snapshot / DriverFactory"
# commit: the card-payment columns (numero do cartao / cvv / mes / ano) are
# dropped from the Cadastro sheet, the Lupa ("magnifying glass" search) sheet
# is repurposed into a simple two-row "item to search" lookup box instead of
# a live formula mirror of Cadastro, and the Home sheet is emptied out.

$wb = $excel.ActiveWorkbook

$cadastro = $wb.Worksheets.Item("Cadastro")
$lupa     = $wb.Worksheets.Item("Lupa")
$home     = $wb.Worksheets.Item("Home")

# ---------------------------------------------------------------------------
# Lupa: turn it into a small "item to search" box.
#  - A1/A2 become plain labels instead of formulas pulled from Cadastro.
#  - B1/B2 (the senha / Test12 mirror) are removed entirely.
#  - C1:F2 (the old card-detail formula mirror) are cleared back to blank,
#    visually-styled cells (underlined, matching the leftover H5 style already
#    used on Cadastro).
# ---------------------------------------------------------------------------
$lupa.Range("B1:B2").Clear()
$lupa.Range("C1:F2").ClearContents()
$lupa.Range("C1:F2").Font.Underline = 2

$lupa.Range("A1").Value = "Item a ser pesquisado"
$lupa.Range("A2").Value = "HP Pavilion 15z Laptop"

$lupa.Columns.Item(1).AutoFit()

$lupa.PageSetup.PaperSize = 9
$lupa.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Cadastro: drop the card columns (K:N -> numero do cartao, cvv, mes, ano)
# and rename the sample person in A2.
# ---------------------------------------------------------------------------
$cadastro.Range("K1:N2").Clear()
$cadastro.Range("A2").Value = "gggggggggg"

# ---------------------------------------------------------------------------
# Home: wipe the cached A2/B2 formula mirror, leaving an empty sheet.
# ---------------------------------------------------------------------------
$home.Cells.Clear()

# Leave Cadastro's own selection parked back on A2 (first edited cell) and
# Home's selection spanning the old A1:B2 block, then finish on Lupa so it
# is the sheet/tab that ends up active when the workbook is saved.
[void]$cadastro.Range("A2").Select()
[void]$home.Range("A1:B2").Select()

$lupa.Activate()
[void]$lupa.Range("A6").Select()
